$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.446.34'
$ws.Range("E2").Value = '  +3.83%  '
$ws.Range("D3").Value = '2.317.26'
$ws.Range("E3").Value = '  +2.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.64%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.80'
$ws.Range("E5").Value = '  +4.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.38'
$ws.Range("E6").Value = '  +4.31%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +2.47%  '
$ws.Range("D9").Value = '2.339.13'
$ws.Range("E9").Value = '  +2.95%  '
$ws.Range("E10").Value = '  +9.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("E12").Value = '  +6.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  +2.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.03'
$ws.Range("E14").Value = '  +4.80%  '
$ws.Range("D15").Value = '2.733.09'
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("D16").Value = '56.877.26'
$ws.Range("E16").Value = '  +4.97%  '
$ws.Range("E17").Value = '  +4.91%  '
$ws.Range("D18").Value = '2.329.00'
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.54'
$ws.Range("E19").Value = '  +3.05%  '
$ws.Range("E20").Value = '  +3.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.50'
$ws.Range("E21").Value = '  +6.39%  '
$ws.Range("E22").Value = '  +4.91%  '
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.17'
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.70'
$ws.Range("E27").Value = '  +5.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.37'
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("E29").Value = '  +11.99%  '
$ws.Range("D30").Value = '0.0₃0739'
$ws.Range("E30").Value = '  +6.99%  '
$ws.Range("E31").Value = '  +5.17%  '
$ws.Range("E32").Value = '  +4.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.35'
$ws.Range("E33").Value = '  +3.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.992'
$ws.Range("E35").Value = '  -0.60%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.27'
$ws.Range("E36").Value = '  +5.16%  '
$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.946'
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  +8.91%  '
$ws.Range("E39").Value = '  +8.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.59'
$ws.Range("E40").Value = '  +4.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.382'
$ws.Range("E41").Value = '  +1.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.71'
$ws.Range("E42").Value = '  +12.96%  '
$ws.Range("E43").Value = '  +7.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '277.95'
$ws.Range("E44").Value = '  +14.91%  '
$ws.Range("E45").Value = '  +7.51%  '
$ws.Range("E46").Value = '  +3.76%  '
$ws.Range("E47").Value = '  +4.07%  '
$ws.Range("E48").Value = '  +3.19%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0216'
$ws.Range("E49").Value = '  +5.94%  '
$ws.Range("B50").Value = 'Polygon'
$ws.Range("C50").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.382'
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.95'
$ws.Range("E51").Value = '  +5.17%  '
